# Update odds values on Sheet1 to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G3"  = 2.05
    "H3"  = 3
    "I3"  = 3.9
    "J3"  = 2.88
    "K3"  = 1.95
    "M3"  = 1.11
    "N3"  = 6.5
    "Q3"  = 1.85
    "R3"  = 2
    "AD3" = 19
    "AE3" = 21
    "AF3" = 41
    "AG3" = 6.5
    "AP3" = 51

    "M4"  = 1.06
    "N4"  = 10
    "S4"  = 2.08
    "T4"  = 1.73

    "Q5"  = 2.1
    "R5"  = 1.78

    "G6"  = 1.8
    "J6"  = 2.5
    "AB6" = 8.5
    "AC6" = 9
    "AD6" = 15
    "AH6" = 6

    "H7"  = 3.5
    "J7"  = 4.33
    "K7"  = 2.25
    "O7"  = 1.25
    "P7"  = 3.75
    "S7"  = 1.83
    "T7"  = 2.03
    "U7"  = 3
    "V7"  = 1.36
    "W7"  = 1.36
    "X7"  = 3
    "Y7"  = 1.67
    "Z7"  = 2.1
    "AG7" = 12
    "AK7" = 8
    "AL7" = 9.5
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
